$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SMX.PotatoHab","bases",11),
    @("SMX.AsteroidPort","bases",11),
    @("SMX.Size0Generator","electrics",7),
    @("SMX.Size1Generator","electrics",10),
    @("SMX.Mk2Generator","electrics",11),
    @("SMX.Size2Generator","electrics",11),
    @("SMX.Radialklaw","isru",9),
    @("SMX.Size0ISRU","isru",9),
    @("SME.MiniKlaw","isru",9),
    @("SMX.Size0OreTank","isru",9),
    @("SMX.0mStackDrill","isru",9),
    @("SMX.Mk2Oretank","isru",10),
    @("SMX.Mk2ISRURefinery","isru",10),
    @("SMX.Size1AtmCondenser","isru",10),
    @("SMX.Size0Driver","isru",10),
    @("SMX.1mInlineDrill","isru",10),
    @("SMX.1mHInlineDrill","isru",10),
    @("SMX.1mInlinePump","isru",10),
    @("SMX.1mStackDrill","isru",10),
    @("SMX.1mStackPump","isru",10),
    @("SMX.InlineDrill","isru",11),
    @("SMX.Mk3ISRURefinery","isru",11),
    @("SMX.Mk3OreTank","isru",11),
    @("SMX.RCSDriver","isru",11),
    @("SMX.Size2AtmCondenser","isru",11),
    @("SMX.Size1Driver","isru",11),
    @("SMX.2mStackDrill","isru",11),
    @("SMX.2mStackPump","isru",11),
    @("SMX.Size3ISRU","isru",12),
    @("SMX.Size3OreTank","isru",12),
    @("SMX.3mStackDrill","isru",12),
    @("SMX.VLandingGear","landing",8),
    @("SMX.VLandingGearL","landing",9),
    @("SMX.ShroudedVLandingGear","landing",9),
    @("SMX.ShroudedVLandingGearL","landing",10),
    @("SMX.RadialPillarLeg","landing",11),
    @("SMX.StackLeg","landing",11),
    @("SMX.DeployableWheel","robotics",8),
    @("SMX.PoddedDeployableWheelS","robotics",9),
    @("SMX.PoddedDeployableWheelL","robotics",10)
)

$row = 2
foreach ($item in $data) {
    $name = $item[0]
    $cat = $item[1]
    $tier = $item[2]

    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $cat
    $ws.Cells.Item($row, 3).Value = $tier
    $ws.Cells.Item($row, 4).Value = "MiningExpansion"
    $ws.Cells.Item($row, 5).Formula = '="@PART[" & A' + $row + ' & "]:AFTER[" & D' + $row + ' & "] //' + "`n" + '{' + "`n" + "`t" + '@TechRequired = " & B' + $row + ' & C' + $row + ' & "' + "`n" + '}"'

    $row++
}

$ws.Columns.Item(4).ColumnWidth = 28.44140625

$ws.Range("A19").Select()
$ws.Range("E41").Select()
